## adding 3-year county data and assoicated organizing; starting to add "education" strata
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variableNames")
$ws.Activate()

# Insert a new row above the current last data row (row 15), pushing the
# existing "F62"/"countyFIPS" row down to row 16, and fill in the new
# "education" / "F35" strata row.
$ws.Rows("15:15").Insert()

$ws.Range("A15").Value = "education"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "F35"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 201
$ws.Range("F15").Value = 201

# Match the new selection left behind on the sheet (D16, below the
# newly-inserted row, same column as the header style column).
$ws.Range("D16").Select() | Out-Null

# Printer/page setup tweak that shipped in the same commit.
$ws.PageSetup.Orientation = 1
